$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-04-01 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-02 Wednesday", 2) | Out-Null

# Update each table cell directly by position to avoid cross-cell text collisions
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "12+16="
$t.Cell(1, 2).Range.Text = "20-17="
$t.Cell(1, 3).Range.Text = "72-48="
$t.Cell(1, 4).Range.Text = "53-30="
$t.Cell(1, 5).Range.Text = "11+62="

$t.Cell(2, 1).Range.Text = "35+32="
$t.Cell(2, 2).Range.Text = "40+50="
$t.Cell(2, 3).Range.Text = "23+61="
$t.Cell(2, 4).Range.Text = "55-25="
$t.Cell(2, 5).Range.Text = "91-57="

$t.Cell(3, 1).Range.Text = "69+26="
$t.Cell(3, 2).Range.Text = "67+30="
$t.Cell(3, 3).Range.Text = "19+49="
$t.Cell(3, 4).Range.Text = "95-94="
$t.Cell(3, 5).Range.Text = "67+31="

$t.Cell(4, 1).Range.Text = "19-0="
$t.Cell(4, 2).Range.Text = "82-32="
$t.Cell(4, 3).Range.Text = "29+39="
$t.Cell(4, 4).Range.Text = "95-36="
$t.Cell(4, 5).Range.Text = "76-36="

$t.Cell(5, 1).Range.Text = "36+28="
$t.Cell(5, 2).Range.Text = "92-71="
$t.Cell(5, 3).Range.Text = "92-10="
$t.Cell(5, 4).Range.Text = "25+56="
$t.Cell(5, 5).Range.Text = "33+14="

$t.Cell(6, 1).Range.Text = "46+48="
$t.Cell(6, 2).Range.Text = "29-15="
$t.Cell(6, 3).Range.Text = "43-19="
$t.Cell(6, 4).Range.Text = "53-53="
$t.Cell(6, 5).Range.Text = "7+33="

$t.Cell(7, 1).Range.Text = "13-12="
$t.Cell(7, 2).Range.Text = "11+17="
$t.Cell(7, 3).Range.Text = "21+60="
$t.Cell(7, 4).Range.Text = "67-59="
$t.Cell(7, 5).Range.Text = "19+51="

$t.Cell(8, 1).Range.Text = "80-78="
$t.Cell(8, 2).Range.Text = "78-7="
$t.Cell(8, 3).Range.Text = "84+9="
$t.Cell(8, 4).Range.Text = "62-52="
$t.Cell(8, 5).Range.Text = "73-44="

$t.Cell(9, 1).Range.Text = "79-31="
$t.Cell(9, 2).Range.Text = "34-9="
$t.Cell(9, 3).Range.Text = "44-22="
$t.Cell(9, 4).Range.Text = "38+48="
$t.Cell(9, 5).Range.Text = "72-35="

$t.Cell(10, 1).Range.Text = "53+37="
$t.Cell(10, 2).Range.Text = "83-15="
$t.Cell(10, 3).Range.Text = "96-63="
$t.Cell(10, 4).Range.Text = "92-51="
$t.Cell(10, 5).Range.Text = "37+0="

$t.Cell(11, 1).Range.Text = "44+37="
$t.Cell(11, 2).Range.Text = "43-30="
$t.Cell(11, 3).Range.Text = "67+4="
$t.Cell(11, 4).Range.Text = "90-77="
$t.Cell(11, 5).Range.Text = "32+52="

$t.Cell(12, 1).Range.Text = "57-37="
$t.Cell(12, 2).Range.Text = "43+12="
$t.Cell(12, 3).Range.Text = "50-17="
$t.Cell(12, 4).Range.Text = "90-69="
$t.Cell(12, 5).Range.Text = "55-49="

$t.Cell(13, 1).Range.Text = "63-46="
$t.Cell(13, 2).Range.Text = "24+58="
$t.Cell(13, 3).Range.Text = "41-10="
$t.Cell(13, 4).Range.Text = "55-42="
$t.Cell(13, 5).Range.Text = "26+41="

$t.Cell(14, 1).Range.Text = "6+58="
$t.Cell(14, 2).Range.Text = "16-10="
$t.Cell(14, 3).Range.Text = "44-33="
$t.Cell(14, 4).Range.Text = "80+0="
$t.Cell(14, 5).Range.Text = "24+59="

$t.Cell(15, 1).Range.Text = "12+43="
$t.Cell(15, 2).Range.Text = "12+66="
$t.Cell(15, 3).Range.Text = "2+4="
$t.Cell(15, 4).Range.Text = "36-31="
$t.Cell(15, 5).Range.Text = "89-3="

$t.Cell(16, 1).Range.Text = "4+82="
$t.Cell(16, 2).Range.Text = "96-38="
$t.Cell(16, 3).Range.Text = "91-87="
$t.Cell(16, 4).Range.Text = "75-39="
$t.Cell(16, 5).Range.Text = "78+5="

$t.Cell(17, 1).Range.Text = "8-7="
$t.Cell(17, 2).Range.Text = "74+14="
$t.Cell(17, 3).Range.Text = "15+42="
$t.Cell(17, 4).Range.Text = "31+35="
$t.Cell(17, 5).Range.Text = "90-41="

$t.Cell(18, 1).Range.Text = "71-65="
$t.Cell(18, 2).Range.Text = "33+29="
$t.Cell(18, 3).Range.Text = "21+42="
$t.Cell(18, 4).Range.Text = "45+54="
$t.Cell(18, 5).Range.Text = "36+24="

$t.Cell(19, 1).Range.Text = "22+65="
$t.Cell(19, 2).Range.Text = "92-62="
$t.Cell(19, 3).Range.Text = "54-52="
$t.Cell(19, 4).Range.Text = "79-17="
$t.Cell(19, 5).Range.Text = "32+20="

$t.Cell(20, 1).Range.Text = "97-61="
$t.Cell(20, 2).Range.Text = "47+33="
$t.Cell(20, 3).Range.Text = "7+92="
$t.Cell(20, 4).Range.Text = "34+55="
$t.Cell(20, 5).Range.Text = "12+64="
